$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2, 1).Value = "Última actualización: 08:21:27"
$ws1.Cells.Item(3, 1).Value = "Total filas: 106"

$ws1.Cells.Item(73, 1).Value = "07:46:15"
$ws1.Cells.Item(73, 2).Value = "08:23"
$ws1.Cells.Item(73, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(73, 4).Value = 37
$ws1.Cells.Item(73, 5).Value = "LP1912"

$ws1.Cells.Item(74, 1).Value = "06:52:23"
$ws1.Cells.Item(74, 2).Value = "08:23"
$ws1.Cells.Item(74, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(74, 4).Value = 91
$ws1.Cells.Item(74, 5).Value = "LP1912"

$ws1.Cells.Item(75, 1).Value = "06:52:23"
$ws1.Cells.Item(75, 2).Value = "08:27"
$ws1.Cells.Item(75, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(75, 4).Value = 95
$ws1.Cells.Item(75, 5).Value = "LP1912"

$ws1.Cells.Item(76, 1).Value = "07:46:15"
$ws1.Cells.Item(76, 2).Value = "08:33"
$ws1.Cells.Item(76, 3).Value = "10_OLMOS"
$ws1.Cells.Item(76, 4).Value = 47
$ws1.Cells.Item(76, 5).Value = "LP1912"

$ws1.Cells.Item(77, 1).Value = "07:46:15"
$ws1.Cells.Item(77, 2).Value = "08:33"
$ws1.Cells.Item(77, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(77, 4).Value = 47
$ws1.Cells.Item(77, 5).Value = "LP1912"

$ws1.Cells.Item(78, 1).Value = "07:46:15"
$ws1.Cells.Item(78, 2).Value = "08:34"
$ws1.Cells.Item(78, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(78, 4).Value = 48
$ws1.Cells.Item(78, 5).Value = "LP1912"

$ws1.Cells.Item(79, 1).Value = "08:21:27"
$ws1.Cells.Item(79, 2).Value = "08:34"
$ws1.Cells.Item(79, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(79, 4).Value = 13
$ws1.Cells.Item(79, 5).Value = "LP1912"

$ws1.Cells.Item(80, 1).Value = "07:59:05"
$ws1.Cells.Item(80, 2).Value = "08:35"
$ws1.Cells.Item(80, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(80, 4).Value = 36
$ws1.Cells.Item(80, 5).Value = "LP1912"

$ws1.Cells.Item(81, 1).Value = "06:52:23"
$ws1.Cells.Item(81, 2).Value = "08:42"
$ws1.Cells.Item(81, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(81, 4).Value = 110
$ws1.Cells.Item(81, 5).Value = "LP1912"

$ws1.Cells.Item(82, 1).Value = "08:21:27"
$ws1.Cells.Item(82, 2).Value = "08:43"
$ws1.Cells.Item(82, 3).Value = "14_ABASTO"
$ws1.Cells.Item(82, 4).Value = 22
$ws1.Cells.Item(82, 5).Value = "LP1912"

$ws1.Cells.Item(83, 1).Value = "07:46:15"
$ws1.Cells.Item(83, 2).Value = "08:44"
$ws1.Cells.Item(83, 3).Value = "14_ABASTO"
$ws1.Cells.Item(83, 4).Value = 58
$ws1.Cells.Item(83, 5).Value = "LP1912"

$ws1.Cells.Item(84, 1).Value = "07:59:05"
$ws1.Cells.Item(84, 2).Value = "08:48"
$ws1.Cells.Item(84, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(84, 4).Value = 49
$ws1.Cells.Item(84, 5).Value = "LP1912"

$ws1.Cells.Item(85, 1).Value = "07:59:05"
$ws1.Cells.Item(85, 2).Value = "08:53"
$ws1.Cells.Item(85, 3).Value = "10_OLMOS"
$ws1.Cells.Item(85, 4).Value = 54
$ws1.Cells.Item(85, 5).Value = "LP1912"

$ws1.Cells.Item(86, 1).Value = "07:46:15"
$ws1.Cells.Item(86, 2).Value = "08:54"
$ws1.Cells.Item(86, 3).Value = "17_ROMERO"
$ws1.Cells.Item(86, 4).Value = 68
$ws1.Cells.Item(86, 5).Value = "LP1912"

$ws1.Cells.Item(87, 1).Value = "08:21:27"
$ws1.Cells.Item(87, 2).Value = "09:01"
$ws1.Cells.Item(87, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(87, 4).Value = 40
$ws1.Cells.Item(87, 5).Value = "LP1912"

$ws1.Cells.Item(88, 1).Value = "08:21:27"
$ws1.Cells.Item(88, 2).Value = "09:01"
$ws1.Cells.Item(88, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(88, 4).Value = 40
$ws1.Cells.Item(88, 5).Value = "LP1912"

$ws1.Cells.Item(89, 1).Value = "07:46:15"
$ws1.Cells.Item(89, 2).Value = "09:02"
$ws1.Cells.Item(89, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(89, 4).Value = 76
$ws1.Cells.Item(89, 5).Value = "LP1912"

$ws1.Cells.Item(90, 1).Value = "08:21:27"
$ws1.Cells.Item(90, 2).Value = "09:03"
$ws1.Cells.Item(90, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(90, 4).Value = 42
$ws1.Cells.Item(90, 5).Value = "LP1912"

$ws1.Cells.Item(91, 1).Value = "07:46:15"
$ws1.Cells.Item(91, 2).Value = "09:04"
$ws1.Cells.Item(91, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(91, 4).Value = 78
$ws1.Cells.Item(91, 5).Value = "LP1912"

$ws1.Cells.Item(92, 1).Value = "08:21:27"
$ws1.Cells.Item(92, 2).Value = "09:09"
$ws1.Cells.Item(92, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(92, 4).Value = 48
$ws1.Cells.Item(92, 5).Value = "LP1912"

$ws1.Cells.Item(93, 1).Value = "08:21:27"
$ws1.Cells.Item(93, 2).Value = "09:10"
$ws1.Cells.Item(93, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(93, 4).Value = 49
$ws1.Cells.Item(93, 5).Value = "LP1912"

$ws1.Cells.Item(94, 1).Value = "07:46:15"
$ws1.Cells.Item(94, 2).Value = "09:11"
$ws1.Cells.Item(94, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(94, 4).Value = 85
$ws1.Cells.Item(94, 5).Value = "LP1912"

$ws1.Cells.Item(95, 1).Value = "08:21:27"
$ws1.Cells.Item(95, 2).Value = "09:16"
$ws1.Cells.Item(95, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(95, 4).Value = 55
$ws1.Cells.Item(95, 5).Value = "LP1912"

$ws1.Cells.Item(96, 1).Value = "07:46:15"
$ws1.Cells.Item(96, 2).Value = "09:17"
$ws1.Cells.Item(96, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(96, 4).Value = 91
$ws1.Cells.Item(96, 5).Value = "LP1912"

$ws1.Cells.Item(97, 1).Value = "07:46:15"
$ws1.Cells.Item(97, 2).Value = "09:21"
$ws1.Cells.Item(97, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(97, 4).Value = 95
$ws1.Cells.Item(97, 5).Value = "LP1912"

$ws1.Cells.Item(98, 1).Value = "08:21:27"
$ws1.Cells.Item(98, 2).Value = "09:22"
$ws1.Cells.Item(98, 3).Value = "17_ROMERO"
$ws1.Cells.Item(98, 4).Value = 61
$ws1.Cells.Item(98, 5).Value = "LP1912"

$ws1.Cells.Item(99, 1).Value = "07:46:15"
$ws1.Cells.Item(99, 2).Value = "09:22"
$ws1.Cells.Item(99, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(99, 4).Value = 96
$ws1.Cells.Item(99, 5).Value = "LP1912"

$ws1.Cells.Item(100, 1).Value = "07:59:05"
$ws1.Cells.Item(100, 2).Value = "09:23"
$ws1.Cells.Item(100, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(100, 4).Value = 84
$ws1.Cells.Item(100, 5).Value = "LP1912"

$ws1.Cells.Item(101, 1).Value = "07:46:15"
$ws1.Cells.Item(101, 2).Value = "09:23"
$ws1.Cells.Item(101, 3).Value = "17_ROMERO"
$ws1.Cells.Item(101, 4).Value = 97
$ws1.Cells.Item(101, 5).Value = "LP1912"

$ws1.Cells.Item(102, 1).Value = "08:21:27"
$ws1.Cells.Item(102, 2).Value = "09:23"
$ws1.Cells.Item(102, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(102, 4).Value = 62
$ws1.Cells.Item(102, 5).Value = "LP1912"

$ws1.Cells.Item(103, 1).Value = "07:46:15"
$ws1.Cells.Item(103, 2).Value = "09:24"
$ws1.Cells.Item(103, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(103, 4).Value = 98
$ws1.Cells.Item(103, 5).Value = "LP1912"

$ws1.Cells.Item(104, 1).Value = "07:46:15"
$ws1.Cells.Item(104, 2).Value = "09:32"
$ws1.Cells.Item(104, 3).Value = "15_ABASTO"
$ws1.Cells.Item(104, 4).Value = 106
$ws1.Cells.Item(104, 5).Value = "LP1912"

$ws1.Cells.Item(105, 1).Value = "07:46:15"
$ws1.Cells.Item(105, 2).Value = "09:33"
$ws1.Cells.Item(105, 3).Value = "10_OLMOS"
$ws1.Cells.Item(105, 4).Value = 107
$ws1.Cells.Item(105, 5).Value = "LP1912"

$ws1.Cells.Item(106, 1).Value = "08:21:27"
$ws1.Cells.Item(106, 2).Value = "09:36"
$ws1.Cells.Item(106, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(106, 4).Value = 75
$ws1.Cells.Item(106, 5).Value = "LP1912"

$ws1.Cells.Item(107, 1).Value = "07:46:15"
$ws1.Cells.Item(107, 2).Value = "09:42"
$ws1.Cells.Item(107, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(107, 4).Value = 116
$ws1.Cells.Item(107, 5).Value = "LP1912"

$ws1.Cells.Item(108, 1).Value = "08:21:27"
$ws1.Cells.Item(108, 2).Value = "09:43"
$ws1.Cells.Item(108, 3).Value = "14_ABASTO"
$ws1.Cells.Item(108, 4).Value = 82
$ws1.Cells.Item(108, 5).Value = "LP1912"

$ws1.Cells.Item(109, 1).Value = "07:59:05"
$ws1.Cells.Item(109, 2).Value = "09:44"
$ws1.Cells.Item(109, 3).Value = "14_ABASTO"
$ws1.Cells.Item(109, 4).Value = 105
$ws1.Cells.Item(109, 5).Value = "LP1912"

$ws1.Cells.Item(110, 1).Value = "07:59:05"
$ws1.Cells.Item(110, 2).Value = "09:52"
$ws1.Cells.Item(110, 3).Value = "15_ABASTO"
$ws1.Cells.Item(110, 4).Value = 113
$ws1.Cells.Item(110, 5).Value = "LP1912"

$ws1.Cells.Item(111, 1).Value = "08:21:27"
$ws1.Cells.Item(111, 2).Value = "10:12"
$ws1.Cells.Item(111, 3).Value = "15_ABASTO"
$ws1.Cells.Item(111, 4).Value = 111
$ws1.Cells.Item(111, 5).Value = "LP1912"

$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = "Última actualización: 08:21:27"
$ws2.Cells.Item(3, 1).Value = "Total filas: 16"

$ws2.Cells.Item(19, 1).Value = "08:21:27"
$ws2.Cells.Item(19, 2).Value = "09:01"
$ws2.Cells.Item(19, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(19, 4).Value = 40
$ws2.Cells.Item(19, 5).Value = "LP1912"

$ws2.Cells.Item(20, 1).Value = "07:46:15"
$ws2.Cells.Item(20, 2).Value = "09:02"
$ws2.Cells.Item(20, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(20, 4).Value = 76
$ws2.Cells.Item(20, 5).Value = "LP1912"

$ws2.Cells.Item(21, 1).Value = "07:46:15"
$ws2.Cells.Item(21, 2).Value = "09:42"
$ws2.Cells.Item(21, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(21, 4).Value = 116
$ws2.Cells.Item(21, 5).Value = "LP1912"

$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = "Última actualización: 08:21:27"
$ws3.Cells.Item(3, 1).Value = "Total filas: 24"

$ws3.Cells.Item(26, 1).Value = "08:21:27"
$ws3.Cells.Item(26, 2).Value = "08:52"
$ws3.Cells.Item(26, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(26, 4).Value = 31
$ws3.Cells.Item(26, 5).Value = "L6173"

$ws3.Cells.Item(27, 1).Value = "07:16:53"
$ws3.Cells.Item(27, 2).Value = "09:08"
$ws3.Cells.Item(27, 3).Value = "215D_LA PLATA"
$ws3.Cells.Item(27, 4).Value = 112
$ws3.Cells.Item(27, 5).Value = "L6203"

$ws3.Cells.Item(28, 1).Value = "07:46:15"
$ws3.Cells.Item(28, 2).Value = "09:09"
$ws3.Cells.Item(28, 3).Value = "215D_LA PLATA"
$ws3.Cells.Item(28, 4).Value = 83
$ws3.Cells.Item(28, 5).Value = "L6203"

$ws3.Cells.Item(29, 1).Value = "08:21:27"
$ws3.Cells.Item(29, 2).Value = "10:03"
$ws3.Cells.Item(29, 3).Value = "215B_LP-P MOR-40 Y 115"
$ws3.Cells.Item(29, 4).Value = 102
$ws3.Cells.Item(29, 5).Value = "L6173"

